$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 24 and 25: Toncoin/WrappedBTC swap order, with updated values
$ws.Cells.Item(24, 2).Value = "WrappedBTC"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "21.778.34"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -1.86%  "

$ws.Cells.Item(25, 2).Value = "Toncoin"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "2.353"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +3.52%  "

# Remaining per-row Price (D) and Volume(1h) (E) updates
$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "21.781.78"
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  -1.85%  "
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "1.542.87"
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  -1.22%  "
$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = "0.9996"
$c.Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  -0.07%  "
$ws.Cells.Item(5, 5).Value = "  +0.00%  "
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "290.31"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -0.10%  "
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "0.3898"
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +3.35%  "
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "0.3188"
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -3.18%  "
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "43.16"
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -1.26%  "
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "0.07204"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -2.39%  "
$ws.Cells.Item(11, 5).Value = "  -6.39%  "
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "0.9997"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -0.08%  "
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "5.644"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -3.42%  "
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "18.65"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -6.96%  "
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "6.619"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -3.87%  "
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "1.547.60"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -0.37%  "
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "0.00001105"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +0.57%  "
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "0.06568"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  -1.02%  "
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "83.18"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -2.84%  "
$ws.Cells.Item(20, 5).Value = "  -0.01%  "
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "6.158"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -4.60%  "
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "15.36"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -5.02%  "
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "10.90"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -7.40%  "
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "2.401"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -5.52%  "
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "144.44"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -4.18%  "
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "18.42"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -3.65%  "
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "4.842"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -1.27%  "
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "1.718.10"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -0.65%  "
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "117.69"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -3.32%  "
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "0.9727"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -13.87%  "
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "5.926"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -2.26%  "
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "0.08208"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +0.08%  "
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "9.001"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -4.09%  "
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "0.06119"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -1.61%  "
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "5.145"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -2.99%  "
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "0.02217"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -4.37%  "
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "0.2048"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -4.62%  "
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "1.187"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -4.41%  "
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "1.424"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -24.01%  "
$ws.Cells.Item(42, 5).Value = "  -0.02%  "
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "10.63"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -4.26%  "
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "0.5793"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -3.66%  "
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "13.13"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -4.49%  "
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "3.742"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -0.74%  "
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "0.5553"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -4.62%  "
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "1.886"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -5.56%  "
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "117.23"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -3.37%  "
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "1.133"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -3.47%  "
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "0.06733"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -3.86%  "
